$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44285
$ws.Range("K2").Value = 'Angeleno'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 806

# Row 3
$ws.Range("D3").Value = 44243
$ws.Range("K3").Value = 'Black Amber'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 806

# Row 4
$ws.Range("D4").Value = 44314
$ws.Range("K4").Value = 'Angeleno'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 806

# Row 5
$ws.Range("D5").Value = 44175
$ws.Range("K5").Value = 'Angeleno'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1194

# Row 6
$ws.Range("D6").Value = 44217
$ws.Range("K6").Value = 'Black Amber'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 16500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 917

# Row 7
$ws.Range("D7").Value = 44174
$ws.Range("K7").Value = 'Angeleno'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1139

# Row 8
$ws.Range("D8").Value = 44278
$ws.Range("K8").Value = 'Angeleno'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 861

# Row 9
$ws.Range("D9").Value = 44245
$ws.Range("K9").Value = 'Black Amber'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = '$/bandeja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 806

# Row 10
$ws.Range("D10").Value = 44229
$ws.Range("K10").Value = 'Fortuna'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 806

# Row 11
$ws.Range("D11").Value = 44169
$ws.Range("K11").Value = 'Angeleno'
$ws.Range("L11").Value = 'Tercera'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 24500
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1361

# Row 12
$ws.Range("D12").Value = 44239
$ws.Range("K12").Value = 'Fortuna'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15500
$ws.Range("Q12").Value = '$/bandeja 18 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 861

# Row 13
$ws.Range("D13").Value = 44238
$ws.Range("K13").Value = 'Black Amber'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 14500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 806

# Row 14
$ws.Range("D14").Value = 44238
$ws.Range("K14").Value = 'Fortuna'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 14500
$ws.Range("Q14").Value = '$/bandeja 18 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 806
